# Apply the movie-ratings swap described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 163 and 164: swap "The Seventh Seal (1957)" and "The Elephant Man (1980)"
$ws.Range("B163").Value = "The Elephant Man"
$ws.Range("C163").Value = "(1980)"
$ws.Range("B164").Value = "The Seventh Seal"
$ws.Range("C164").Value = "(1957)"

# Rows 248-251: rotate entries and introduce "Soorarai Pottru (2020)"
# Force "96" to be stored as text (it looks numeric otherwise).
$ws.Range("B248").Value = "'96"
$ws.Range("C248").Value = "(2018)"
$ws.Range("B249").Value = "Fanny and Alexander"
$ws.Range("C249").Value = "(1982)"
$ws.Range("B250").Value = "Hera Pheri"
$ws.Range("C250").Value = "(2000)"
$ws.Range("B251").Value = "Soorarai Pottru"
$ws.Range("C251").Value = "(2020)"
